# Levee removal scripts work: rename the "Stream Recharge" budget-term label to
# "Stream Losses" on the owhm_wb_dict sheet, while preserving the original label
# in a new "name_old" column (D) for reference.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("owhm_wb_dict")

# New header for column D
$ws.Range("D2").Value = "name_old"

# Copy the current "name" column (C) values into the new "name_old" column (D)
# for every data row, before any renaming happens.
for ($r = 3; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $ws.Cells.Item($r, 3).Text
}

# Rename the SFR_IN entry from "Stream Recharge" to "Stream Losses"
$ws.Range("C4").Value = "Stream Losses"

# Let column widths adjust to fit the new/duplicated content
$ws.Columns("C:D").AutoFit()

# Leave the cursor where the author left off
$ws.Range("E5").Select() | Out-Null
